# BrickPicker layout.xlsx update
# Simplifies the example table on Tabelle1: collapses the "Whole Color" /
# "Whole Category" / "Green Bricks" example rows into a smaller "All"-based
# table, drops the per-row numeric "ID" column values (keeping a renamed
# "Brick ID" header), and removes the last example row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- header row 21: "ID" column becomes "Brick ID" --------------------
$ws.Range("D21").Value = "Brick ID"
$ws.Range("E21").Value = "Color"

# --- row 23: was "All / Bricks (3) / 3 / Unknown (0, 51, 178)" --------
#     becomes "All / Bricks (3) / <blank> / All" (default colored text)
$ws.Range("B23").Value = "All"
$ws.Range("E23").Value = "All"
$ws.Range("E23").ClearFormats()
$ws.Range("D23").ClearContents()

# --- row 24: was "Whole Color / Whole Category / 4 / Green (0,255,0)" -
#     becomes "All / All / <blank> / Green (0, 255, 0)"
$ws.Range("B24").Value = "All"
$ws.Range("C24").Value = "All"
$ws.Range("D24").ClearContents()

# --- row 25: was "Whole Category / Bricks (3) / 3 / Unknown (...)" ----
#     becomes "All / Bricks (3) / <blank> / Green (0, 255, 0)"
$ws.Range("B25").Value = "All"
$ws.Range("E25").Value = "Green (0, 255, 0)"
$ws.Range("E25").Font.Color = 5287936
$ws.Range("D25").ClearContents()
$ws.Range("G25").Value = "#all green parts in cat. 3"

# --- row 26 ("Green Bricks" example) is removed entirely --------------
#     deleting the row shifts everything below it up by one (old row 30
#     becomes row 29).
$ws.Rows("26:26").Delete()

# the formatted-but-empty placeholder cell that was E30 (now E29) keeps
# its dark-blue font
$ws.Range("E29").Font.Color = 11678464

# --- column D ("Brick ID") needs a sensible width now it has data -----
$ws.Columns("D:D").AutoFit()

# --- selection cosmetics, matching the authored file -------------------
$ws.Range("E31").Select()
